$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Chickcen"
$ws.Range("A10").Value = "Rat"

$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 2

$ws.Range("C9").Value = "Worry"
$ws.Range("C10").Value = "Clever"

$ws.Range("I11").Select()
